$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was "Inhib_E" (2,0,0,0) -> becomes "Inhib_A" (2,0,2,0)
$ws.Range("A2").Value = "Inhib_A"
$ws.Range("D2").Value = 2

# Row 3 was "Inhib_B" (2,0,0,0) -> becomes "Inhib_B/E" (2,0,0,0) - label only change
$ws.Range("A3").Value = "Inhib_B/E"

# Row 4 "Inhib_C/D" (2,2,2,2) unchanged

# Old row 5 ("Inhib_A", 2,0,2,0) was a duplicate of the now-corrected row 2; delete it
$ws.Range("A5:F5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Update the active selection to reflect the saved view state
$ws.Range("A6").Select()
